$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 242, shifting the existing data (old rows
# 242-270) down to rows 245-273. This reproduces the dimension growing
# from A1:T270 to A1:T273 and keeps all of the "downstream" rows' values
# (and the D-column date style) identical to before, as required by the
# diff.
$ws.Rows("242:244").Insert()

# Populate the three newly inserted rows with this week's price report
# (week of 2023-08-16, serial 45154) for the three apple varieties.

# Row 242: Fuji royal, Calibre 100
$ws.Cells.Item(242, 1).Value = 1
$ws.Cells.Item(242, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(242, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(242, 4).Value = 45154
$ws.Cells.Item(242, 5).Value = 15
$ws.Cells.Item(242, 6).Value = "Fruta"
$ws.Cells.Item(242, 7).Value = 100104
$ws.Cells.Item(242, 8).Value = "Frutos de pepita"
$ws.Cells.Item(242, 9).Value = 100104002
$ws.Cells.Item(242, 10).Value = "Manzana"
$ws.Cells.Item(242, 11).Value = "Fuji royal"
$ws.Cells.Item(242, 12).Value = "Calibre 100"
$ws.Cells.Item(242, 13).Value = 270
$ws.Cells.Item(242, 14).Value = 24000
$ws.Cells.Item(242, 15).Value = 25000
$ws.Cells.Item(242, 16).Value = 24500
$ws.Cells.Item(242, 17).Value = "`$/caja 18 kilos embalada"
$ws.Cells.Item(242, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(242, 19).Value = 1361
$ws.Cells.Item(242, 20).Value = 18

# Row 243: Granny Smith, Calibre 90
$ws.Cells.Item(243, 1).Value = 1
$ws.Cells.Item(243, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(243, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(243, 4).Value = 45154
$ws.Cells.Item(243, 5).Value = 15
$ws.Cells.Item(243, 6).Value = "Fruta"
$ws.Cells.Item(243, 7).Value = 100104
$ws.Cells.Item(243, 8).Value = "Frutos de pepita"
$ws.Cells.Item(243, 9).Value = 100104002
$ws.Cells.Item(243, 10).Value = "Manzana"
$ws.Cells.Item(243, 11).Value = "Granny Smith"
$ws.Cells.Item(243, 12).Value = "Calibre 90"
$ws.Cells.Item(243, 13).Value = 300
$ws.Cells.Item(243, 14).Value = 24000
$ws.Cells.Item(243, 15).Value = 25000
$ws.Cells.Item(243, 16).Value = 24500
$ws.Cells.Item(243, 17).Value = "`$/caja 18 kilos embalada"
$ws.Cells.Item(243, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(243, 19).Value = 1361
$ws.Cells.Item(243, 20).Value = 18

# Row 244: Royal Gala, Calibre 100
$ws.Cells.Item(244, 1).Value = 1
$ws.Cells.Item(244, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(244, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(244, 4).Value = 45154
$ws.Cells.Item(244, 5).Value = 15
$ws.Cells.Item(244, 6).Value = "Fruta"
$ws.Cells.Item(244, 7).Value = 100104
$ws.Cells.Item(244, 8).Value = "Frutos de pepita"
$ws.Cells.Item(244, 9).Value = 100104002
$ws.Cells.Item(244, 10).Value = "Manzana"
$ws.Cells.Item(244, 11).Value = "Royal Gala"
$ws.Cells.Item(244, 12).Value = "Calibre 100"
$ws.Cells.Item(244, 13).Value = 300
$ws.Cells.Item(244, 14).Value = 24000
$ws.Cells.Item(244, 15).Value = 25000
$ws.Cells.Item(244, 16).Value = 24500
$ws.Cells.Item(244, 17).Value = "`$/caja 18 kilos embalada"
$ws.Cells.Item(244, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(244, 19).Value = 1361
$ws.Cells.Item(244, 20).Value = 18
